# feat: add 2022-Q4 data
#
#  - insert a new "2022-Q4" worksheet right after the "总计" summary sheet;
#    it carries the per-fund holdings detail for 2022-Q4 (the existing
#    "2022-Q3" detail sheet simply shifts down to the 3rd tab)
#  - update the "总计" summary sheet: the old Q3 totals row becomes row 3,
#    and row 2 is overwritten with the new Q4 totals

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)

# --- insert the new "2022-Q4" sheet right after "总计" ---------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# --- "总计" summary sheet ---------------------------------------------------
# push the old Q3 totals row down to row 3 (copy A2's number style onto the
# new A3 cell so it matches the rest of the column) ...
$summary.Range("A3").Value = 1
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.19

# ... then overwrite row 2 with the new Q4 totals
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.26

# --- "2022-Q4" sheet: fund-level holding detail -----------------------------
# header row + column A share the same bold/bordered style as the "总计"
# sheet's header row / first data column, so copy that formatting over
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# numeric-looking fund codes / ratios must stay TEXT (leading zeros in fund
# codes, fixed decimal display, etc.) instead of being auto-coerced to
# numbers: stamp "@" (text) NumberFormat right before each assignment, then
# drop back to the default "Normal" style so no stray formatting lingers on
# the cell afterwards (matches the source data, which carries no style here)
$textCols = "B2","B3","B4","B5","D2","D3","D4","D5","E2","E3","E4","E5","F2","F3","F4","F5","G2","G3","G4","G5"
foreach ($addr in $textCols) {
    $q4.Range($addr).NumberFormat = "@"
}

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "001411"
$q4.Range("C2").Value = "诺安创新驱动灵活配置混合A"
$q4.Range("D2").Value = "7.22"
$q4.Range("E2").Value = "91.15"
$q4.Range("F2").Value = "2.18"
$q4.Range("G2").Value = "0.1574"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "002051"
$q4.Range("C3").Value = "诺安创新驱动灵活配置混合C"
$q4.Range("D3").Value = "4.60"
$q4.Range("E3").Value = "91.15"
$q4.Range("F3").Value = "2.18"
$q4.Range("G3").Value = "0.1003"
$q4.Range("H3").Value = 8

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "011150"
$q4.Range("C4").Value = "创金合信ESG责任投资股票C"
$q4.Range("D4").Value = "0.12"
$q4.Range("E4").Value = "90.04"
$q4.Range("F4").Value = "2.51"
$q4.Range("G4").Value = "0.0030"
$q4.Range("H4").Value = 4

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "011149"
$q4.Range("C5").Value = "创金合信ESG责任投资股票A"
$q4.Range("D5").Value = "0.10"
$q4.Range("E5").Value = "90.04"
$q4.Range("F5").Value = "2.51"
$q4.Range("G5").Value = "0.0025"
$q4.Range("H5").Value = 4

foreach ($addr in $textCols) {
    $q4.Range($addr).Style = "Normal"
}

# keep "2022-Q3" the selected/active tab, same as before the edit
# (re-fetch by name: sheet indices shifted once "2022-Q4" was inserted)
$wb.Worksheets.Item("2022-Q3").Activate()
